$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6348862741258756
$ws.Range("C2").Value = 0.1660329477771043
$ws.Range("D2").Value = 0.07741023265785429
$ws.Range("E2").Value = 0.1237860245306699
$ws.Range("G2").Value = 0.002522169166861987
$ws.Range("K2").Value = 0.301206131509133
$ws.Range("L2").Value = 0.1936269566045041
$ws.Range("M2").Value = 0.1761975111647409
$ws.Range("O2").Value = 6.019464854163488

$ws.Range("B3").Value = 0.6040193108883614
$ws.Range("C3").Value = 0.1651346122813457
$ws.Range("D3").Value = 0.07038428588825241
$ws.Range("E3").Value = 0.1240181189210272
$ws.Range("G3").Value = 0.002525249183918784
$ws.Range("K3").Value = 0.2726903643780219
$ws.Range("L3").Value = 0.1911695842255696
$ws.Range("M3").Value = 0.1703838107122486
$ws.Range("O3").Value = 6.019389523647959

$ws.Range("B4").Value = 0.5853860983739878
$ws.Range("C4").Value = 0.1645701085249023
$ws.Range("D4").Value = 0.0661056833492637
$ws.Range("E4").Value = 0.1242013713758965
$ws.Range("G4").Value = 0.002527241856007538
$ws.Range("K4").Value = 0.2552809379992595
$ws.Range("L4").Value = 0.1897537074724909
$ws.Range("M4").Value = 0.1669053173570987
$ws.Range("O4").Value = 6.021859274864568

$ws.Range("B5").Value = 0.577873595309228
$ws.Range("C5").Value = 0.1643368222911796
$ws.Range("D5").Value = 0.06437099703272509
$ws.Range("E5").Value = 0.1242863093233364
$ws.Range("G5").Value = 0.002528079493852908
$ws.Range("K5").Value = 0.2482117171233114
$ws.Range("L5").Value = 0.1892001463597097
$ws.Range("M5").Value = 0.1655108034750867
$ws.Range("O5").Value = 6.023498478760786

$ws.Range("B6").Value = 0.5766310350457786
$ws.Range("C6").Value = 0.1642978893097755
$ws.Range("D6").Value = 0.0640834896182696
$ws.Range("E6").Value = 0.1243010333991403
$ws.Range("G6").Value = 0.002528220132099059
$ws.Range("K6").Value = 0.2470394122663748
$ws.Range("L6").Value = 0.1891096438000446
$ws.Range("M6").Value = 0.1652806371000537
$ws.Range("O6").Value = 6.023808887242296

$ws.Range("B7").Value = 0.5852844549058887
$ws.Range("C7").Value = 0.1645669754857515
$ws.Range("D7").Value = 0.06608225282731439
$ws.Range("E7").Value = 0.1242024753115842
$ws.Range("G7").Value = 0.002527253048869439
$ws.Range("K7").Value = 0.255185497364522
$ws.Range("L7").Value = 0.1897461470664226
$ws.Range("M7").Value = 0.166886417218155
$ws.Range("O7").Value = 6.021878819615324

$ws.Range("B8").Value = 0.6241773235777259
$ws.Range("C8").Value = 0.165725885239425
$ws.Range("D8").Value = 0.07498033847961949
$ws.Range("E8").Value = 0.1238576029477798
$ws.Range("G8").Value = 0.00252321013124213
$ws.Range("K8").Value = 0.2913534255243349
$ws.Range("L8").Value = 0.1927603863055154
$ws.Range("M8").Value = 0.1741740774633946
$ws.Range("O8").Value = 6.018916724202853

$ws.Range("B9").Value = 0.7029656906039747
$ws.Range("C9").Value = 0.1678959277226113
$ws.Range("D9").Value = 0.09271147801162272
$ws.Range("E9").Value = 0.1235040041859019
$ws.Range("G9").Value = 0.002516084025555325
$ws.Range("K9").Value = 0.3630588624062625
$ws.Range("L9").Value = 0.1994074638843415
$ws.Range("M9").Value = 0.1891856784037707
$ws.Range("O9").Value = 6.033073731351237

$ws.Range("B10").Value = 0.7623764140513174
$ws.Range("C10").Value = 0.1694277934688557
$ws.Range("D10").Value = 0.1059138734697029
$ws.Range("E10").Value = 0.1234402296680841
$ws.Range("G10").Value = 0.002511332456939115
$ws.Range("K10").Value = 0.4162111265986539
$ws.Range("L10").Value = 0.2047387190339407
$ws.Range("M10").Value = 0.2006517892731594
$ws.Range("O10").Value = 6.055658907306395

$ws.Range("B11").Value = 0.7897331222930575
$ws.Range("C11").Value = 0.1701111397274744
$ws.Range("D11").Value = 0.111958860438591
$ws.Range("E11").Value = 0.1234536371230739
$ws.Range("G11").Value = 0.002509274880822738
$ws.Range("K11").Value = 0.4404928243598079
$ws.Range("L11").Value = 0.2072610590793573
$ws.Range("M11").Value = 0.2059625653375647
$ws.Range("O11").Value = 6.068582458309947

$ws.Range("B12").Value = 0.8001396112309749
$ws.Range("C12").Value = 0.1703679620925627
$ws.Range("D12").Value = 0.1142536028529832
$ws.Range("E12").Value = 0.1234648005123553
$ws.Range("G12").Value = 0.002508510596247519
$ws.Range("K12").Value = 0.4497022171227911
$ws.Range("L12").Value = 0.2082301380198714
$ws.Range("M12").Value = 0.2079871874704651
$ws.Range("O12").Value = 6.073857363544619

$ws.Range("B13").Value = 0.7978962997111125
$ws.Range("C13").Value = 0.1703127374089277
$ws.Range("D13").Value = 0.1137591381969969
$ws.Range("E13").Value = 0.1234621258024617
$ws.Range("G13").Value = 0.002508674538272536
$ws.Range("K13").Value = 0.447718172011804
$ws.Range("L13").Value = 0.2080208109629069
$ws.Range("M13").Value = 0.2075505477094168
$ws.Range("O13").Value = 6.072704370785516

$ws.Range("B14").Value = 0.7905883289238318
$ws.Range("C14").Value = 0.1701323076803973
$ws.Range("D14").Value = 0.1121475372639367
$ws.Range("E14").Value = 0.1234544336410188
$ws.Range("G14").Value = 0.002509211704907772
$ws.Range("K14").Value = 0.441250198905351
$ws.Range("L14").Value = 0.2073405069845933
$ws.Range("M14").Value = 0.2061288612035241
$ws.Range("O14").Value = 6.06900879109935

$ws.Range("B15").Value = 0.7861181067025598
$ws.Range("C15").Value = 0.1700215357226398
$ws.Range("D15").Value = 0.1111611194727971
$ws.Range("E15").Value = 0.1234505141828919
$ws.Range("G15").Value = 0.002509542670284427
$ws.Range("K15").Value = 0.4372902505531044
$ws.Range("L15").Value = 0.2069256131324408
$ws.Range("M15").Value = 0.2052597988523246
$ws.Range("O15").Value = 6.066794765601003

$ws.Range("B16").Value = 0.7605952024580915
$ws.Range("C16").Value = 0.1693828630639445
$ws.Range("D16").Value = 0.1055196077967651
$ws.Range("E16").Value = 0.1234402058797492
$ws.Range("G16").Value = 0.002511469009688011
$ws.Range("K16").Value = 0.4146262972755039
$ws.Range("L16").Value = 0.2045758289118424
$ws.Range("M16").Value = 0.2003066177379367
$ws.Range("O16").Value = 6.054867640688286

$ws.Range("B17").Value = 0.745022066661079
$ws.Range("C17").Value = 0.1689875958745688
$ws.Range("D17").Value = 0.1020687558784772
$ws.Range("E17").Value = 0.1234447369295175
$ws.Range("G17").Value = 0.002512677326540434
$ws.Range("K17").Value = 0.4007487304887718
$ws.Range("L17").Value = 0.2031591578986536
$ws.Range("M17").Value = 0.1972922187151624
$ws.Range("O17").Value = 6.048229357096005

$ws.Range("B18").Value = 0.7360959416679975
$ws.Range("C18").Value = 0.1687589774984417
$ws.Range("D18").Value = 0.1000876018116656
$ws.Range("E18").Value = 0.1234513365330958
$ws.Range("G18").Value = 0.002513382105542288
$ws.Range("K18").Value = 0.3927763831613049
$ws.Range("L18").Value = 0.2023534706618193
$ws.Range("M18").Value = 0.1955673426759219
$ws.Range("O18").Value = 6.044660566774638

$ws.Range("B19").Value = 0.7330790688299089
$ws.Range("C19").Value = 0.1686813530396449
$ws.Range("D19").Value = 0.0994174497083975
$ws.Range("E19").Value = 0.1234542573384392
$ws.Range("G19").Value = 0.002513622414769111
$ws.Range("K19").Value = 0.390078751785353
$ws.Range("L19").Value = 0.2020822506892017
$ws.Range("M19").Value = 0.1949848649784087
$ws.Range("O19").Value = 6.043495064552275

$ws.Range("B20").Value = 0.7466766347099281
$ws.Range("C20").Value = 0.1690298043276215
$ws.Range("D20").Value = 0.1024357235520199
$ws.Range("E20").Value = 0.1234438413810626
$ws.Range("G20").Value = 0.002512547686704066
$ws.Range("K20").Value = 0.402225023313548
$ws.Range("L20").Value = 0.2033090188411251
$ws.Range("M20").Value = 0.1976121833784745
$ws.Range("O20").Value = 6.048910203977385

$ws.Range("B21").Value = 0.7927335822564316
$ws.Range("C21").Value = 0.1701853571083376
$ws.Range("D21").Value = 0.1126207504813692
$ws.Range("E21").Value = 0.1234565279428388
$ws.Range("G21").Value = 0.002509053522893202
$ws.Range("K21").Value = 0.4431496098050616
$ws.Range("L21").Value = 0.2075399512161624
$ws.Range("M21").Value = 0.2065460780075767
$ws.Range("O21").Value = 6.070083930795533

$ws.Range("B22").Value = 0.8231087834411142
$ws.Range("C22").Value = 0.1709292330611945
$ws.Range("D22").Value = 0.1193101143600757
$ws.Range("E22").Value = 0.123500288586893
$ws.Range("G22").Value = 0.002506856554452101
$ws.Range("K22").Value = 0.4699802299360272
$ws.Range("L22").Value = 0.2103862470979294
$ws.Range("M22").Value = 0.2124638154954255
$ws.Range("O22").Value = 6.086143184190803

$ws.Range("B23").Value = 0.8068720058808196
$ws.Range("C23").Value = 0.1705332516878642
$ws.Range("D23").Value = 0.1157368660174711
$ws.Range("E23").Value = 0.123473691719127
$ws.Range("G23").Value = 0.002508021210877174
$ws.Range("K23").Value = 0.4556526301723807
$ws.Range("L23").Value = 0.2088597155523217
$ws.Range("M23").Value = 0.2092982138269477
$ws.Range("O23").Value = 6.077368820014442

$ws.Range("B24").Value = 0.7459285201190369
$ws.Range("C24").Value = 0.1690107261667251
$ws.Range("D24").Value = 0.1022698087892024
$ws.Range("E24").Value = 0.1234442338147783
$ws.Range("G24").Value = 0.002512606265412649
$ws.Range("K24").Value = 0.4015575724599216
$ws.Range("L24").Value = 0.2032412393708398
$ws.Range("M24").Value = 0.1974675019866652
$ws.Range("O24").Value = 6.048601621674067

$ws.Range("B25").Value = 0.6813827679351618
$ws.Range("C25").Value = 0.167319848666736
$ws.Range("D25").Value = 0.08788418979175106
$ws.Range("E25").Value = 0.1235651887489973
$ws.Range("G25").Value = 0.002517926476916976
$ws.Range("K25").Value = 0.343577675305113
$ws.Range("L25").Value = 0.1975305600891701
$ws.Range("M25").Value = 0.1850477308785372
$ws.Range("O25").Value = 6.027105144737334
